$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with the new job posting details
$ws.Range("B2").Value = "Communications Analyst"

# Force the posting date to stay a literal text string (e.g. "01/16/2026")
# instead of being auto-converted to a date serial number, then strip the
# formatting back off so no stray cell style is left behind.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "01/16/2026"
$ws.Range("D2").ClearFormats()

$ws.Range("E2").Formula = '=HYPERLINK("https://estm.fa.em2.oraclecloud.com/hcmUI/CandidateExperience/en/sites/CX_1/job/31299/?location=India&locationId=300000000440677&locationLevel=country&mode=location", "Apply")'

# Remove row 3 (the old "Data Analytics and Strategic Insights Analyst" posting)
$ws.Rows("3").Delete()
